$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 41717
$ws1.Range("G2").Value = "已售罄"
$ws1.Range("F5").Value = 9437
$ws1.Range("F6").Value = 199
$ws1.Range("F7").Value = 842
$ws1.Range("F8").Value = 890
$ws1.Range("F9").Value = 718
$ws1.Range("F10").Value = 209
$ws1.Range("F12").Value = 292
$ws1.Range("F13").Value = 886
$ws1.Range("F15").Value = 122
$ws1.Range("F16").Value = 724
$ws1.Range("F18").Value = 1383
$ws1.Range("F20").Value = 646
$ws1.Range("F21").Value = 694
$ws1.Range("F23").Value = 678
$ws1.Range("F24").Value = 721
$ws1.Range("F27").Value = 60
$ws1.Range("F28").Value = 493
$ws1.Range("F29").Value = 514
$ws1.Range("F30").Value = 46
$ws1.Range("F31").Value = 232
$ws1.Range("F32").Value = 923
$ws1.Range("F35").Value = 90
$ws1.Range("F38").Value = 378
$ws1.Range("F39").Value = 1240
$ws1.Range("F40").Value = 287
$ws1.Range("F42").Value = 1222
$ws1.Range("F43").Value = 370
$ws1.Range("F48").Value = 44
$ws1.Range("F49").Value = 63

# --- Sheet: 本地生活 ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 2023
$ws3.Range("F3").Value = 515
$ws3.Range("F4").Value = 383

# --- Sheet: 全部类型 ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 2023
$ws4.Range("F3").Value = 515
$ws4.Range("F4").Value = 41717
$ws4.Range("G4").Value = "已售罄"
$ws4.Range("F10").Value = 9437
$ws4.Range("F11").Value = 199
$ws4.Range("F12").Value = 842
$ws4.Range("F13").Value = 842
$ws4.Range("F15").Value = 383
$ws4.Range("F16").Value = 890
$ws4.Range("F18").Value = 209
$ws4.Range("F19").Value = 292
$ws4.Range("F20").Value = 886
$ws4.Range("F24").Value = 724
$ws4.Range("F26").Value = 1383
$ws4.Range("F27").Value = 646
$ws4.Range("F28").Value = 694
$ws4.Range("F30").Value = 678
$ws4.Range("F31").Value = 721
$ws4.Range("F33").Value = 60
$ws4.Range("F34").Value = 493
$ws4.Range("F35").Value = 46
$ws4.Range("F36").Value = 232
$ws4.Range("F37").Value = 923
$ws4.Range("F40").Value = 90
$ws4.Range("F42").Value = 378
$ws4.Range("F43").Value = 1222
$ws4.Range("F44").Value = 370
$ws4.Range("F48").Value = 44
$ws4.Range("F50").Value = 63
